# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 00:35"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1834977
$ws.Range("C4").Value = 18157
$ws.Range("D4").Value = 541326
$ws.Range("E4").Value = 1187505
$ws.Range("G4").Value = 589
$ws.Range("H4").Value = 106146

# Japon (row 45)
$ws.Range("B45").Value = 16851
$ws.Range("C45").Value = 47
$ws.Range("D45").Value = 14459
$ws.Range("E45").Value = 1501
$ws.Range("G45").Value = 5
$ws.Range("H45").Value = 891

# Nigeria (row 56)
$ws.Range("B56").Value = 10162
$ws.Range("C56").Value = 307
$ws.Range("D56").Value = 3007
$ws.Range("E56").Value = 6868
$ws.Range("G56").Value = 14
$ws.Range("H56").Value = 287

# Malaui (row 155)
$ws.Range("B155").Value = 284
$ws.Range("C155").Value = 5
$ws.Range("E155").Value = 238

# Zimbabue (row 162)
$ws.Range("B162").Value = 178
$ws.Range("C162").Value = 4
$ws.Range("E162").Value = 145
